$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run-boundary at the edges of a Range without altering its
# final formatting. We do this by toggling Bold off-then-on (or on-then-off,
# depending on its current value) so the range ends up re-serialized as its
# own run(s), but its Font.Bold value is restored to what it originally was.
# ---------------------------------------------------------------------------
function Split-Range($rng) {
    $orig = $rng.Font.Bold
    if ($orig) {
        $rng.Font.Bold = 0
        $rng.Font.Bold = 1
    } else {
        $rng.Font.Bold = 1
        $rng.Font.Bold = 0
    }
}

# ===========================================================================
# 1) "Web" + " Framework" (2 italic runs) -> merge into a single run
#    "Web Framework" (italic). A same-text Find/Replace across the two runs
#    consolidates them.
# ===========================================================================
$d.Content.Find.Execute("Web Framework", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Web Framework", 2) | Out-Null

# ===========================================================================
# 2) Remove the hidden "_GoBack" bookmark after "Instalação e configuração
#    do Banco de dados."
# ===========================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ===========================================================================
# 3) "a página [1]." -> "a página " + "ManualConfiguracaoBD, no Wiki."
#    (split into two runs)
# ===========================================================================
$p7 = $d.Paragraphs(7)
$rng7 = $p7.Range
$rng7.Find.ClearFormatting()
$rng7.Find.Execute("[1].", $true, $false, $false, $false, $false, `
    $true, 1, $false, "ManualConfiguracaoBD, no Wiki.", 2) | Out-Null

# ===========================================================================
# 4) Grails paragraph: split the single run into many runs with new wording.
#    First replace the trailing content, then split at the required offsets.
# ===========================================================================
$p9 = $d.Paragraphs(9)
$rng9 = $p9.Range
$rng9.Find.ClearFormatting()
$rng9.Find.Execute( `
    "O arquivo [2] é um zip contendo o core do framework e é necessário para o projeto.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Em  [1] podemos achar o zip contendo o core do framework que é necessário para o projeto.", `
    2) | Out-Null

# ===========================================================================
# 5) Tomcat paragraph: merge "...versão " + "[3]" + " a partir da 6.X." into
#    a single run (same text, Find/Replace consolidates the runs).
# ===========================================================================
$p11 = $d.Paragraphs(11)
$rng11 = $p11.Range
$rng11.Find.ClearFormatting()
$tomcatText = "Para rodar o Neo SI é necessário ter o servidor Web Tomcat instalado. Pode-se baixar qualquer versão [3] a partir da 6.X."
$rng11.Find.Execute($tomcatText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $tomcatText, 2) | Out-Null

# ===========================================================================
# 6) WAR generation paragraph: replace tail wording, to be split afterwards.
# ===========================================================================
$p13 = $d.Paragraphs(13)
$rng13 = $p13.Range
$rng13.Find.ClearFormatting()
$rng13.Find.Execute( `
    "grails war. Será gerado um arquivo .WAR na pasta target do projeto.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "tomcat deploy. O plugin do tomcat já está integrado ao projeto, então será gerado um arquivo .WAR na pasta target do projeto. ", `
    2) | Out-Null

# ===========================================================================
# 7) Remove the stale <w:lastRenderedPageBreak/> cached before the "Re" run
#    of "Referências Externas", while preserving the existing "Re" /
#    "ferências Externas" run split and all formatting.
# ===========================================================================
$rngRe = $d.Content
$rngRe.Find.ClearFormatting()
$rngRe.Find.Execute("Referências Externas", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$reStart = $rngRe.Start
$firstChar = $d.Range($reStart, $reStart + 1)
$firstChar.Text = "X"
$firstChar2 = $d.Range($reStart, $reStart + 1)
$firstChar2.Text = "R"

# ===========================================================================
# Now perform all the run-splits (must happen *after* all text replacements,
# since Word tends to re-merge adjacent same-formatted runs whenever a
# Find/Replace touches them).
# ===========================================================================

# --- split paragraph 7: "a página " | "ManualConfiguracaoBD, no Wiki."
$p7b = $d.Paragraphs(7)
$pStart7 = $p7b.Range.Start
$cut7 = 92 + 9  # length of 1st (unmodified) run + length of "a página "
Split-Range($d.Range($pStart7, $pStart7 + $cut7))

# --- split paragraph 9 at multiple offsets
$p9b = $d.Paragraphs(9)
$pStart9 = $p9b.Range.Start
$pEnd9 = $p9b.Range.End
$cuts9 = @(64, 67, 69, 70, 71, 88, 121, 123)
$bounds9 = @(0) + $cuts9 + @(($pEnd9 - $pStart9) - 1)
for ($i = 0; $i -lt $bounds9.Length - 1; $i++) {
    $s = $pStart9 + $bounds9[$i]
    $e = $pStart9 + $bounds9[$i + 1]
    Split-Range($d.Range($s, $e))
}

# --- split paragraph 13 at multiple offsets
$p13b = $d.Paragraphs(13)
$pStart13 = $p13b.Range.Start
$pEnd13 = $p13b.Range.End
$cuts13 = @(117, 130, 131, 132, 149, 188, 215)
$bounds13 = @(0) + $cuts13 + @(($pEnd13 - $pStart13) - 1)
for ($i = 0; $i -lt $bounds13.Length - 1; $i++) {
    $s = $pStart13 + $bounds13[$i]
    $e = $pStart13 + $bounds13[$i + 1]
    Split-Range($d.Range($s, $e))
}

# --- restore/ensure "Re" / "ferências Externas" split survives the
#     lastRenderedPageBreak removal trick above.
$rngRe2 = $d.Content
$rngRe2.Find.ClearFormatting()
$rngRe2.Find.Execute("Referências Externas", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
Split-Range($d.Range($rngRe2.Start, $rngRe2.Start + 2))

Write-Host "done"
